$d = $word.ActiveDocument

# Find `text` literally, searching forward from character offset `from`.
# Returns the matched Range (Start/End give its position).
function Find-Range([int]$from, [string]$text) {
    $docEnd = $d.Content.End
    $r = $d.Range($from, $docEnd)
    $ok = $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Could not find text: $text"
    }
    return $r
}

# Find `text` literally and replace it with `replacement`, using an
# explicit Delete + InsertBefore so that any markup anchored purely to
# the old run boundaries (proofErr spell/gram markers, bookmarks, ...)
# is actually dropped rather than silently preserved by a same-length
# in place ".Text =" edit.
function Replace-Text([int]$from, [string]$text, [string]$replacement) {
    $r = Find-Range $from $text
    $start = $r.Start
    $r.Delete()
    $ins = $d.Range($start, $start)
    $ins.InsertBefore($replacement)
    return $d.Range($start, $start + $replacement.Length)
}

# ---------------------------------------------------------------
# Change 1: split "...the line numbering." into "...the l" | "ine
# numbering." and drop a _GoBack bookmark at the split point (this is
# the first occurrence of this sentence, in the Referee #1 preamble).
# ---------------------------------------------------------------
$r1 = Find-Range 0 "Disabling comments/tracked changes will change the line numbering."
$splitPos = $r1.Start + ("Disabling comments/tracked changes will change the l").Length
$d.Bookmarks.Add("_GoBack", $d.Range($splitPos, $splitPos)) | Out-Null

# ---------------------------------------------------------------
# Change 2: "...new stacked bar graph (html plots in Data S2 and S3). We..."
#         -> "...new stacked bar graph (Data S1). We..."
# ---------------------------------------------------------------
$r2 = Replace-Text 0 "html plots in Data S2 and S3" "Data S1"

# ---------------------------------------------------------------
# Change 3: "...Site 1 and Site 2 (html plots in Data S2 and S3)." ->
#           "...Site 1 and Site 2 (Data S1). Data S1 file is best viewed in xlsx format."
# ---------------------------------------------------------------
$r3 = Replace-Text $r2.End "html plots in Data S2 and S3" "Data S1"
$r3b = Replace-Text $r3.End ")." "). Data S1 file is best viewed in xlsx format."

# ---------------------------------------------------------------
# Change 4: rewrite the "In response to this comment, we added several
# supplementary data files: ..." paragraph.
# ---------------------------------------------------------------
$oldP4 = "supplementary data files:. Data S1 contains the OTU tables of both Site 1 and Site 2 with OTU representative sequences, taxonomy, and abundances across the timeline. Data S2 and S3 include the raw OTU information, and interactive stacked bar plots of the community across time timelines for both sites. Data S4 contains details about MAG taxonomy, statistics, and abundance.  "
$newP4 = "supplementary data files: Data S1 contains the OTU tables of both Site 1 and Site 2 with OTU representative sequences, taxonomy, abundances across the timeline, and stacked bar plots of the community across time. Data S2 contains details about MAG taxonomy, statistics, and abundance. Data S1 and S2 files are best viewed in xlsx format."
$r4 = Replace-Text $r3b.End $oldP4 $newP4

# ---------------------------------------------------------------
# Change 5: merge the run split around "its" (drops the spell-check
# proofErr wrapper and the stray extra space run) - no visible text change.
# ---------------------------------------------------------------
$r5 = Replace-Text $r4.End "to walk the reader through the figure and its " "to walk the reader through the figure and its "

# ---------------------------------------------------------------
# Change 6: merge "(see sup. files" + ")" into a single run and drop the
# _GoBack bookmark that used to sit here (it moved to change 1 above).
# ---------------------------------------------------------------
$r6 = Replace-Text $r5.End " (see sup. files)" " (see sup. files)"

Write-Output "done"
